$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.657073
$ws.Range("H2").Value = 1.971219
$ws.Range("I2").Value = 0.576753533868729
$ws.Range("J2").Value = 0.576753533868729
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 4.143401148758
$ws.Range("R2").Value = 37.290610338822
$ws.Range("S2").Value = 0.00782164600206361
$ws.Range("T2").Value = 0.007821646002063608

# Row 3
$ws.Range("G3").Value = 0.657073
$ws.Range("H3").Value = 1.971219
$ws.Range("I3").Value = 0.576753533868729
$ws.Range("J3").Value = 0.576753533868729
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 119.9366401058104
$ws.Range("R3").Value = 1079.429760952293
$ws.Range("S3").Value = 0.2264086695698661
$ws.Range("T3").Value = 0.226408669569866

# Row 4
$ws.Range("G4").Value = 0.657073
$ws.Range("H4").Value = 1.971219
$ws.Range("I4").Value = 0.576753533868729
$ws.Range("J4").Value = 0.576753533868729
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 83.708730137689
$ws.Range("R4").Value = 753.3785712392011
$ws.Range("S4").Value = 0.1580199529112786
$ws.Range("T4").Value = 0.1580199529112786

# Row 5
$ws.Range("G5").Value = 0.657073
$ws.Range("H5").Value = 1.971219
$ws.Range("I5").Value = 0.576753533868729
$ws.Range("J5").Value = 0.576753533868729
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 12.76298682809733
$ws.Range("R5").Value = 114.866881452876
$ws.Range("S5").Value = 0.0240931450550719
$ws.Range("T5").Value = 0.02409314505507189

# Row 6
$ws.Range("G6").Value = 0.657073
$ws.Range("H6").Value = 1.971219
$ws.Range("I6").Value = 0.576753533868729
$ws.Range("J6").Value = 0.576753533868729
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 84.97488593503668
$ws.Range("R6").Value = 764.77397341533
$ws.Range("S6").Value = 0.1604101203304489
$ws.Range("T6").Value = 0.1604101203304489

# Row 7
$ws.Range("G7").Value = 0.4821883333333334
$ws.Range("H7").Value = 1.446565
$ws.Range("I7").Value = 0.423246466131271
$ws.Range("J7").Value = 0.423246466131271
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 3.040605372996667
$ws.Range("R7").Value = 27.36544835697
$ws.Range("S7").Value = 0.005739859117112379
$ws.Range("T7").Value = 0.005739859117112378

# Row 8
$ws.Range("G8").Value = 0.4821883333333334
$ws.Range("H8").Value = 1.446565
$ws.Range("I8").Value = 0.423246466131271
$ws.Range("J8").Value = 0.423246466131271
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("Q8").Value = 88.01464768483946
$ws.Range("R8").Value = 792.131829163555
$ws.Range("S8").Value = 0.1661483869099949
$ws.Range("T8").Value = 0.1661483869099949

# Row 9
$ws.Range("G9").Value = 0.4821883333333334
$ws.Range("H9").Value = 1.446565
$ws.Range("I9").Value = 0.423246466131271
$ws.Range("J9").Value = 0.423246466131271
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 61.42905441334834
$ws.Range("R9").Value = 552.8614897201351
$ws.Range("S9").Value = 0.1159618150916279
$ws.Range("T9").Value = 0.1159618150916279

# Row 10
$ws.Range("G10").Value = 0.4821883333333334
$ws.Range("H10").Value = 1.446565
$ws.Range("I10").Value = 0.423246466131271
$ws.Range("J10").Value = 0.423246466131271
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 9.366026829584445
$ws.Range("R10").Value = 84.29424146626
$ws.Range("S10").Value = 0.01768058261237847
$ws.Range("T10").Value = 0.01768058261237847

# Row 11
$ws.Range("G11").Value = 0.4821883333333334
$ws.Range("H11").Value = 1.446565
$ws.Range("I11").Value = 0.423246466131271
$ws.Range("J11").Value = 0.423246466131271
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 62.35821381217224
$ws.Range("R11").Value = 561.2239243095501
$ws.Range("S11").Value = 0.1177158224001574
$ws.Range("T11").Value = 0.1177158224001573
